$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New leaderboard rows to append (Customer Name, Salesperson, Prospect, Last Invoice Date(blank), Customer Number)
$newRows = @(
    @{ Row = 11; A = "YOUNG'S";                     B = "Larsen, Rick J"; C = "040"; E = "0008325" },
    @{ Row = 12; A = "MARMA";                        B = "Larsen, Rick J"; C = "040"; E = "0008326" },
    @{ Row = 13; A = "BLOSSOM BRIDGE CHILD CARE";     B = "Larsen, Rick J"; C = "040"; E = "0008327" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Match existing row height/formatting used throughout the sheet
    $ws.Rows.Item($rowNum).RowHeight = 13.05

    # Column A - Customer Name
    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellA.NumberFormat = "@"
    $cellA.VerticalAlignment = -4160
    $cellA.Value = $r.A

    # Column B - Salesperson
    $cellB = $ws.Cells.Item($rowNum, 2)
    $cellB.NumberFormat = "@"
    $cellB.VerticalAlignment = -4160
    $cellB.Value = $r.B

    # Column C - Prospect
    $cellC = $ws.Cells.Item($rowNum, 3)
    $cellC.NumberFormat = "@"
    $cellC.VerticalAlignment = -4160
    $cellC.Value = $r.C

    # Column D - Last Invoice Date (left blank, matches style of other blank rows)
    $cellD = $ws.Cells.Item($rowNum, 4)
    $cellD.VerticalAlignment = -4160

    # Column E - Customer Number
    $cellE = $ws.Cells.Item($rowNum, 5)
    $cellE.NumberFormat = "@"
    $cellE.VerticalAlignment = -4160
    $cellE.Value = $r.E

    # Column F - present as an empty placeholder cell like the rest of the sheet
    $ws.Cells.Item($rowNum, 6).Font.Name = "Arial"
}
